# Add a new "DoubleDQN" run row into the Sheet1 log table.
# Before: rows 2-5 hold DoubleDQN, DoubleDQN, DuellingDQN, DuellingDQN.
# After:  a new DoubleDQN row is inserted at row 4 (pushing the two
#         DuellingDQN rows down to 5 and 6), growing the table to A1:K6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 4, shifting existing rows 4-5 down to 5-6.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new run's data.
$ws.Range("A4").Value = "DoubleDQN"
$ws.Range("B4").Value = 10000
$ws.Range("C4").Value = 32
$ws.Range("D4").Value = 1000000
$ws.Range("E4").Value = 0.1
$ws.Range("F4").Value = 150000
$ws.Range("G4").Value = "YES"
$ws.Range("H4").Value = "NO"
$ws.Range("I4").Value = "running"
$ws.Range("K4").Value = "hlaptop"

# Match the saved cursor/selection position recorded in the file.
$ws.Range("D11").Select()
